# Edit applied to the IES axioms worksheet:
# 1. Shorten the "State" class Comment in C1613 so it no longer includes the
#    two "Note:" paragraphs about whole-life Elements / temporal parts.
# 2. Remove the axiom row that states
#    "State SubClassOf http://ies.data.gov.uk/ontology/ies4#Element"
#    (originally row 1618). Deleting the row shifts every following row up
#    by one, which matches the new sheet dimension of A1:C1930.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: shorten the Comment text on row 1613 (State's Comment axiom)
$ws.Range("C1613").Value = "Comment: A temporal state of an Element"

# Step 2: delete the row containing "SubClassOf: .../ies4#Element" for State
$ws.Rows(1618).Delete()
